$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new rows store numeric-looking values as TEXT (matching the
# existing rows' t="str" cell type) instead of being auto-coerced to numbers.
$ws.Range("A5:K7").NumberFormat = "@"

# Row 5 — duplicate of the "Abu Dhabi" match (same data as row 3)
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " November 06 2020"
$ws.Range("C5").Value = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$ws.Range("D5").Value = "Royal Challengers Bangalore"
$ws.Range("E5").Value = "Sunrisers Hyderabad"
$ws.Range("F5").Value = "Mohammed Siraj "
$ws.Range("G5").Value = "10"
$ws.Range("H5").Value = "7"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "142.85"

# Row 6 — duplicate of the "Sharjah" match (same data as row 2)
$ws.Range("A6").Value = " Sharjah"
$ws.Range("B6").Value = " October 31 2020"
$ws.Range("C6").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D6").Value = "Royal Challengers Bangalore"
$ws.Range("E6").Value = "Sunrisers Hyderabad"
$ws.Range("F6").Value = "Mohammed Siraj "
$ws.Range("G6").Value = "2"
$ws.Range("H6").Value = "3"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "66.66"

# Row 7 — duplicate of the "Dubai (DSC)" match (same data as row 4)
$ws.Range("A7").Value = " Dubai (DSC)"
$ws.Range("B7").Value = " October 05 2020"
$ws.Range("C7").Value = "Capitals won by 59 runs"
$ws.Range("D7").Value = "Royal Challengers Bangalore"
$ws.Range("E7").Value = "Delhi Capitals"
$ws.Range("F7").Value = "Mohammed Siraj "
$ws.Range("G7").Value = "5"
$ws.Range("H7").Value = "4"
$ws.Range("I7").Value = "1"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "125.00"

# Restore the default/"Normal" style so the new cells don't pick up a custom
# number-format style index (keeps styles.xml -> cell s="" attribute aligned
# with the rest of the sheet, which all use the default style).
$ws.Range("A5:K7").Style = "Normal"
